$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values are updated in-place following an automatic electricity price refresh.
# Column A holds the day (date serial); it advances by one day.
$ws.Range("A2").Value = 45923

$ws.Range("B2").Value = 79.11
$ws.Range("C2").Value = 61
$ws.Range("D2").Value = 85.28
$ws.Range("E2").Value = 76.95
$ws.Range("F2").Value = 79.11
$ws.Range("G2").Value = 83.7
$ws.Range("H2").Value = 85.28
$ws.Range("I2").Value = 86.56
$ws.Range("J2").Value = 85.8
$ws.Range("K2").Value = 39.9
$ws.Range("L2").Value = 3.82
$ws.Range("M2").Value = 0.65
$ws.Range("N2").Value = 0.01
$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 4.31
$ws.Range("Q2").Value = 7.9
$ws.Range("R2").Value = 4.31
$ws.Range("S2").Value = 2.2
$ws.Range("T2").Value = 26.47
$ws.Range("U2").Value = 86.31999999999999
$ws.Range("V2").Value = 105.01
$ws.Range("W2").Value = 98.5
$ws.Range("X2").Value = 76.02
$ws.Range("Y2").Value = 75.84
$ws.Range("Z2").Value = 52.27

# AA2 (Slot_4h_max) unchanged: "20h-24h"
$ws.Range("AB2").Value = 88.84
# AC2 (Slot_2h_frist) unchanged: "20h-22h"
$ws.Range("AD2").Value = 101.76
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 85.92
$ws.Range("AG2").Value = "9h-18h"
